$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# PREPARATION (F2): "Kode Transaksi" updated from 090 to 998
$newPrep = "Username : 30711;`nPassword : bni1234;`nRole : 09 - Penyelia Settlement;`nKode Transaksi : 998;`nNama Jenis Transaksi : Saldo Awal Top Up Edit"
$ws.Range("F2").Value = $newPrep

# KODE_JENIS_TRANSAKSI (M2): was the text "090", now the number 998.
# Setting .Value directly resets the cell's style to the default, so grab
# the formatting from a neighbouring cell that shares M2's original style
# (K2, style index 7 - quotePrefix/left+vcenter alignment) and restore it
# after writing the new value.
$ws.Range("M2").Value = 998
$ws.Range("K2").Copy() | Out-Null
$ws.Range("M2").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Update the sheet selection to match the saved view (J2 active cell)
$ws.Range("J2").Select() | Out-Null
